# NYPD CompStat weekly refresh: bump the report volume/number and the
# covered week dates, then replace the crime-stat table (rows 14-30,
# columns C:N) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 32   Number  51" -> "...Number  52" ---
$ws.Range("A8").Value = "Volume 32   Number  52"

# --- Header: reporting week "12/15/2025 ... 12/21/2025" -> "12/22/2025 ... 12/28/2025" ---
$ws.Range("C9").Value = "Report Covering the Week  12/22/2025  Through  12/28/2025"

# --- Row 14 (Murder) ---
$ws.Range("N14").Value = -85.185185185185

# --- Row 15 (Rape) ---
$ws.Range("C15").Value = 4
$ws.Range("F15").Value = 8
$ws.Range("H15").Value = 166.666666666667
$ws.Range("I15").Value = 38
$ws.Range("K15").Value = -15.555555555555
$ws.Range("L15").Value = -2.564102564102
$ws.Range("M15").Value = 31.034482758620
$ws.Range("N15").Value = -62

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 6.25
$ws.Range("I16").Value = 280
$ws.Range("J16").Value = 254
$ws.Range("K16").Value = 10.236220472440
$ws.Range("L16").Value = 10.236220472440
$ws.Range("M16").Value = -32.692307692307
$ws.Range("N16").Value = -88.120492151039

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 40
$ws.Range("G17").Value = 65
$ws.Range("H17").Value = -38.461538461538
$ws.Range("I17").Value = 740
$ws.Range("J17").Value = 850
$ws.Range("K17").Value = -12.941176470588
$ws.Range("L17").Value = 12.804878048780
$ws.Range("M17").Value = 63.716814159292
$ws.Range("N17").Value = -39.690301548492

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 45.454545454545
$ws.Range("I18").Value = 185
$ws.Range("J18").Value = 191
$ws.Range("K18").Value = -3.141361256544
$ws.Range("L18").Value = -2.631578947368
$ws.Range("M18").Value = -55.421686746988
$ws.Range("N18").Value = -91.223908918406

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 26.190476190476
$ws.Range("I19").Value = 686
$ws.Range("J19").Value = 579
$ws.Range("K19").Value = 18.480138169257
$ws.Range("L19").Value = 4.255319148936
$ws.Range("M19").Value = 28.464419475655
$ws.Range("N19").Value = -11.711711711711

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 13
$ws.Range("H20").Value = 62.5
$ws.Range("I20").Value = 199
$ws.Range("J20").Value = 191
$ws.Range("K20").Value = 4.188481675392
$ws.Range("L20").Value = -3.864734299516
$ws.Range("M20").Value = -26.022304832713
$ws.Range("N20").Value = -90.158259149357

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 35
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 147
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = 0.684931506849
$ws.Range("I21").Value = 2136
$ws.Range("J21").Value = 2120
$ws.Range("K21").Value = 0.754716981132
$ws.Range("L21").Value = 5.690252350321
$ws.Range("M21").Value = -0.233535730966
$ws.Range("N21").Value = -75.292076344707

# --- Row 22 (Transit) --- (some cells flip from "N/A" text to real numbers)
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("H22").Value = 0
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 0

# --- Row 23 (Housing) ---
$ws.Range("G23").Value = 3
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = -10

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 31.818181818181
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = -1.834862385321
$ws.Range("I24").Value = 1168
$ws.Range("J24").Value = 1212
$ws.Range("K24").Value = -3.630363036303
$ws.Range("L24").Value = 2.816901408450
$ws.Range("M24").Value = 25.187566988210

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 14
$ws.Range("H25").Value = 27.272727272727
$ws.Range("I25").Value = 209
$ws.Range("J25").Value = 189
$ws.Range("K25").Value = 10.582010582010
$ws.Range("L25").Value = 10.582010582010

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 80
$ws.Range("F26").Value = 70
$ws.Range("G26").Value = 52
$ws.Range("H26").Value = 34.615384615384
$ws.Range("I26").Value = 891
$ws.Range("J26").Value = 950
$ws.Range("K26").Value = -6.210526315789
$ws.Range("L26").Value = 16.166883963494
$ws.Range("M26").Value = 1.365187713310

# --- Row 27 (UCR Rape*) ---
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = 166.666666666667
$ws.Range("I27").Value = 43
$ws.Range("K27").Value = -28.333333333333
$ws.Range("L27").Value = -14

# --- Row 28 (Other Sex Crimes) --- (C28 flips from "N/A" text to a real number)
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 86
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = 22.857142857142
$ws.Range("L28").Value = 19.444444444444

# --- Row 29 (Shooting Vic.) ---
$ws.Range("L29").Value = -11.904761904761
$ws.Range("M29").Value = -54.320987654321
$ws.Range("N29").Value = -83.333333333333

# --- Row 30 (Shooting Inc.) ---
$ws.Range("L30").Value = -17.142857142857
$ws.Range("M30").Value = -58.571428571428
$ws.Range("N30").Value = -85.572139303482
